$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells keep their literal text representation
# (avoids Excel auto-converting numeric-looking strings to numbers,
# which would drop formatting like trailing zeros).
$textCells = @(
    "D2", "D3", "D5", "D6", "D7", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.626.02"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "3.467.69"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "596.80"
$ws.Range("E5").Value = "  -3.30%  "
$ws.Range("D6").Value = "146.74"
$ws.Range("E6").Value = "  -5.11%  "
$ws.Range("D7").Value = "3.467.15"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").Value = "7.85"
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("D12").Value = "0.420"
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("D13").Value = "0.0000211"
$ws.Range("E13").Value = "  -4.45%  "
$ws.Range("D14").Value = "4.051.86"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "31.02"
$ws.Range("E15").Value = "  -6.81%  "
$ws.Range("D16").Value = "3.467.79"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "66.629.40"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").Value = "10.20"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "15.19"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").Value = "432.06"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("E23").Value = "  -6.21%  "
$ws.Range("D24").Value = "79.17"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.604.01"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  -11.83%  "
$ws.Range("D28").Value = "9.67"
$ws.Range("E28").Value = "  -7.82%  "
$ws.Range("D29").Value = "8.18"
$ws.Range("E29").Value = "  -11.74%  "
$ws.Range("D30").Value = "2.46"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "1.58"
$ws.Range("E31").Value = "  -8.28%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  -4.44%  "
$ws.Range("D34").Value = "25.33"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "3.458.89"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.83"
$ws.Range("E37").Value = "  -8.28%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "7.90"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "170.60"
$ws.Range("E41").Value = "  -4.97%  "
$ws.Range("D42").Value = "0.0879"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "5.36"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "2.04"
$ws.Range("E44").Value = "  -14.14%  "
$ws.Range("D45").Value = "0.893"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "45.97"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "28.52"
$ws.Range("E47").Value = "  -9.07%  "
$ws.Range("E48").Value = "  -9.44%  "
$ws.Range("D49").Value = "7.42"
$ws.Range("E49").Value = "  -4.71%  "
$ws.Range("E50").Value = "  -10.07%  "
$ws.Range("D51").Value = "0.960"
$ws.Range("E51").Value = "  -5.57%  "

Write-Output "Applied 95 cell updates"
